# Update cryptos list with latest prices / volume(1h) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to keep a numeric-looking string as TEXT (matches the
    # source data, which stores every Price/Volume figure as an inline
    # string, never a real number).
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "71.634.11"
$ws.Range("E2").Value = "  +3.22%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "4.011.76"
$ws.Range("E3").Value = "  +1.98%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.12%  "

# Row 5 - BNB
Set-TextValue "D5" "528.89"
$ws.Range("E5").Value = "  +2.57%  "

# Row 6 - Solana
Set-TextValue "D6" "148.74"
$ws.Range("E6").Value = "  +1.98%  "

# Row 7 - XRP
Set-TextValue "D7" "0.628"
$ws.Range("E7").Value = "  +1.08%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.31%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.58%  "

# Row 11 - ShibaInu
Set-TextValue "D11" "0.0000345"
$ws.Range("E11").Value = "  +0.66%  "

# Row 12 - Avalanche
Set-TextValue "D12" "44.70"
$ws.Range("E12").Value = "  +3.74%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  +3.57%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.650.58"
$ws.Range("E14").Value = "  +2.13%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "4.023.76"
$ws.Range("E15").Value = "  +2.40%  "

# Row 16 - Chainlink
Set-TextValue "D16" "21.40"
$ws.Range("E16").Value = "  +8.29%  "

# Row 17 - Uniswap
Set-TextValue "D17" "14.33"
$ws.Range("E17").Value = "  +1.38%  "

# Row 18 - Polygon
$ws.Range("E18").Value = "  +0.09%  "

# Row 19 - TRON
Set-TextValue "D19" "0.133"
$ws.Range("E19").Value = "  -1.68%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "71.624.12"
$ws.Range("E20").Value = "  +3.29%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "443.03"
$ws.Range("E21").Value = "  +2.44%  "

# Row 22 - ImmutableX
Set-TextValue "D22" "3.59"
$ws.Range("E22").Value = "  +5.38%  "

# Row 23 - Litecoin
Set-TextValue "D23" "93.83"
$ws.Range("E23").Value = "  +6.38%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" "14.43"

# Row 25 - RenderToken
Set-TextValue "D25" "12.40"
$ws.Range("E25").Value = "  +4.69%  "

# Row 26 - PancakeSwap
Set-TextValue "D26" "4.12"
$ws.Range("E26").Value = "  +5.71%  "

# Row 27 - Filecoin
Set-TextValue "D27" "11.00"
$ws.Range("E27").Value = "  -0.86%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "37.25"
$ws.Range("E28").Value = "  +1.42%  "

# Row 29 - Bittensor
Set-TextValue "D29" "705.62"
$ws.Range("E29").Value = "  +0.50%  "

# Row 30 - Cosmos
Set-TextValue "D30" "13.66"

# Row 31 - Hedera
Set-TextValue "D31" "0.130"
$ws.Range("E31").Value = "  +1.87%  "

# Rows 32/33 - Toncoin and NEARProtocol swap ranking order
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D32" "7.05"
$ws.Range("E32").Value = "  +18.50%  "

$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D33" "2.90"
$ws.Range("E33").Value = "  +1.04%  "

# Row 34 - OKB
Set-TextValue "D34" "68.22"
$ws.Range("E34").Value = "  +2.54%  "

# Row 35 - PEPE
$ws.Range("D35").Value = "0.0$([char]0x2083)0905"
$ws.Range("E35").Value = "  +2.94%  "

# Row 36 - TheGraph
$ws.Range("E36").Value = "  +1.31%  "

# Row 37 - InjectiveProtocol
Set-TextValue "D37" "41.08"
$ws.Range("E37").Value = "  +1.97%  "

# Row 38 - ThetaToken
$ws.Range("E38").Value = "  +21.42%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +2.37%  "

# Row 40 - Dai
$ws.Range("E40").Value = "  +0.32%  "

# Row 41 - VeChain
$ws.Range("E41").Value = "  +2.78%  "

# Row 42 - FirstDigitalUSD
Set-TextValue "D42" "1.00"
$ws.Range("E42").Value = "  +0.00%  "

# Row 43 - Fetch.AI
Set-TextValue "D43" "2.89"
$ws.Range("E43").Value = "  +2.55%  "

# Row 44 - WEMIXToken
Set-TextValue "D44" "3.14"
$ws.Range("E44").Value = "  +1.40%  "

# Row 45 - ApeXProtocol
Set-TextValue "D45" "3.55"
$ws.Range("E45").Value = "  +5.60%  "

# Row 46 - Stacks
Set-TextValue "D46" "3.24"
$ws.Range("E46").Value = "  +9.69%  "

# Row 48 - FLOKI
$ws.Range("E48").Value = "  +22.00%  "

# Row 49 - THORChain
Set-TextValue "D49" "9.32"
$ws.Range("E49").Value = "  +6.89%  "

# Row 50 - LidoDAOToken
$ws.Range("E50").Value = "  +0.87%  "

# Row 51 - BabyDogeCoin
$ws.Range("D51").Value = "0.0$([char]0x2086)0346"
$ws.Range("E51").Value = "  -3.70%  "
